# Add prolificid in rank to use in binary.
# The re_rank score (column E) was recomputed, which reshuffles the
# row order for the female/realeffort ranking. Columns B (prolificid),
# C (name) and F (race) move together with each worker as rows are
# re-sorted by the new E value; A (positional index) and G (rank
# position, 1..12) stay fixed per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Giana (unchanged position, only score refreshed)
$ws.Cells.Item(2, 5).Value = 7.240540192629654

# Row 3 - Jewel moves into this slot (was Colleen)
$ws.Cells.Item(3, 2).Value = 19
$ws.Cells.Item(3, 3).Value = "Jewel"
$ws.Cells.Item(3, 5).Value = 6.378978103426058
$ws.Cells.Item(3, 6).Value = "Black or African American"

# Row 4 - Colleen moves into this slot (was Jewel)
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "Colleen"
$ws.Cells.Item(4, 5).Value = 6.143455313863114
$ws.Cells.Item(4, 6).Value = "White"

# Row 5 - Annes (unchanged position, only score refreshed)
$ws.Cells.Item(5, 5).Value = 5.419772607443591

# Row 6 - Tina (unchanged position, only score refreshed)
$ws.Cells.Item(6, 5).Value = 5.312796240675778

# Row 7 - Nansi (unchanged position, only score refreshed)
$ws.Cells.Item(7, 5).Value = 4.071991992584385

# Row 8 - Khushi (unchanged position, only score refreshed)
$ws.Cells.Item(8, 5).Value = 1.233832614214271

# Row 9 - Lori (unchanged position, only score refreshed)
$ws.Cells.Item(9, 5).Value = 1.002782814522061

# Row 10 - Shaniek (unchanged position, only score refreshed)
$ws.Cells.Item(10, 5).Value = 0.4746561773749075

# Row 11 - Bri moves into this slot (was Kellie)
$ws.Cells.Item(11, 2).Value = 21
$ws.Cells.Item(11, 3).Value = "Bri"
$ws.Cells.Item(11, 5).Value = 0.4664235049697223
$ws.Cells.Item(11, 6).Value = "Black or African American"

# Row 12 - Kellie moves into this slot (was Shadaisia)
$ws.Cells.Item(12, 2).Value = 32
$ws.Cells.Item(12, 3).Value = "Kellie"
$ws.Cells.Item(12, 5).Value = 0.2827895313987391
$ws.Cells.Item(12, 6).Value = "White"

# Row 13 - Shadaisia moves into this slot (was Bri)
$ws.Cells.Item(13, 2).Value = 30
$ws.Cells.Item(13, 3).Value = "Shadaisia"
$ws.Cells.Item(13, 5).Value = 0.27386664857579
